$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.410.43"
$ws.Range("E2").Value = "  +7.19%  "
$ws.Range("D3").Value = "2.378.82"
$ws.Range("E3").Value = "  +4.34%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "112.68"
$ws.Range("E5").Value = "  +9.61%  "
$ws.Range("D6").Value = "317.43"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +4.56%  "
$ws.Range("D10").Value = "42.68"
$ws.Range("E10").Value = "  +11.39%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").Value = "8.69"
$ws.Range("E12").Value = "  +6.24%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("D15").Value = "15.86"
$ws.Range("E15").Value = "  +5.68%  "
$ws.Range("D16").Value = "2.742.59"
$ws.Range("E16").Value = "  +4.58%  "
$ws.Range("D17").Value = "2.377.32"
$ws.Range("E17").Value = "  +4.50%  "
$ws.Range("D18").Value = "45.307.68"
$ws.Range("E18").Value = "  +7.09%  "
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("D22").Value = "74.81"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Value = "3.56"
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("D24").Value = "269.56"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  +9.60%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "7.61"
$ws.Range("E27").Value = "  +7.66%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.76%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "39.31"
$ws.Range("E30").Value = "  +9.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("D32").Value = "0.0955"
$ws.Range("E32").Value = "  +12.78%  "
$ws.Range("D33").Value = "171.17"
$ws.Range("E33").Value = "  +4.21%  "
$ws.Range("E34").Value = "  +15.72%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  +8.37%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.93"
$ws.Range("E37").Value = "  +10.43%  "
$ws.Range("E38").Value = "  +12.07%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "4.01"
$ws.Range("E39").Value = "  +11.24%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0364"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("E41").Value = "  +11.04%  "
$ws.Range("D42").Value = "105.25"
$ws.Range("E42").Value = "  +6.51%  "
$ws.Range("D43").Value = "0.239"
$ws.Range("E43").Value = "  +7.00%  "
$ws.Range("D44").Value = "71.42"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  +11.01%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "5.78"
$ws.Range("E47").Value = "  +13.22%  "
$ws.Range("D48").Value = "116.65"
$ws.Range("E48").Value = "  +6.37%  "
$ws.Range("E49").Value = "  +20.26%  "
$ws.Range("E50").Value = "  +8.42%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "78.75"
$ws.Range("E51").Value = "  +2.85%  "

Write-Host "Updated cryptos list"
